$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 411.29166
$ws.Cells.Item(107, 9).Value = 251.6
$ws.Cells.Item(107, 10).Value = 677.44446
$ws.Cells.Item(107, 11).Value = 251.6
$ws.Cells.Item(107, 12).Value = 677.44446
$ws.Cells.Item(107, 13).Value = 1668.4
$ws.Cells.Item(107, 14).Value = -4517.44446

$ws.Cells.Item(113, 8).Value = 3955.2666
$ws.Cells.Item(113, 9).Value = 4401.2856
$ws.Cells.Item(113, 10).Value = 3565
$ws.Cells.Item(113, 11).Value = 4401.2856
$ws.Cells.Item(113, 12).Value = 3565
$ws.Cells.Item(113, 13).Value = -1147.2856
$ws.Cells.Item(113, 14).Value = -10073

$ws.Cells.Item(127, 8).Value = 1394
$ws.Cells.Item(127, 9).Value = 586.7143
$ws.Cells.Item(127, 10).Value = 1726.4117
$ws.Cells.Item(127, 11).Value = 1760.1429
$ws.Cells.Item(127, 12).Value = 5179.2351
$ws.Cells.Item(127, 13).Value = 3199.8571
$ws.Cells.Item(127, 14).Value = -15099.2351

$ws.Cells.Item(137, 8).Value = 1112.8909
$ws.Cells.Item(137, 9).Value = 1077.2291
$ws.Cells.Item(137, 10).Value = 1357.4286
$ws.Cells.Item(137, 11).Value = 3231.6873
$ws.Cells.Item(137, 12).Value = 4072.2858
$ws.Cells.Item(137, 13).Value = -681.6873000000001
$ws.Cells.Item(137, 14).Value = -9172.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(34, 8).Value = 20000
$ws.Cells.Item(34, 10).Value = 20000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 14).Value = -20542

$ws.Cells.Item(61, 8).Value = 1966.1212
$ws.Cells.Item(61, 9).Value = 1799.36
$ws.Cells.Item(61, 10).Value = 2487.25
$ws.Cells.Item(61, 11).Value = 1799.36
$ws.Cells.Item(61, 12).Value = 2487.25
$ws.Cells.Item(61, 13).Value = -1587.36
$ws.Cells.Item(61, 14).Value = -2911.25

$ws.Cells.Item(74, 8).Value = 936.70215
$ws.Cells.Item(74, 9).Value = 816.6842
$ws.Cells.Item(74, 10).Value = 1443.4445
$ws.Cells.Item(74, 11).Value = 816.6842
$ws.Cells.Item(74, 12).Value = 1443.4445
$ws.Cells.Item(74, 13).Value = 57.31579999999997
$ws.Cells.Item(74, 14).Value = -3191.4445

$ws.Cells.Item(77, 8).Value = 936.70215
$ws.Cells.Item(77, 9).Value = 816.6842
$ws.Cells.Item(77, 10).Value = 1443.4445
$ws.Cells.Item(77, 11).Value = 4083.421
$ws.Cells.Item(77, 12).Value = 7217.2225
$ws.Cells.Item(77, 13).Value = 284.5789999999997
$ws.Cells.Item(77, 14).Value = -15953.2225

$ws.Cells.Item(122, 8).Value = 2446.6
$ws.Cells.Item(122, 9).Value = 2395.7693
$ws.Cells.Item(122, 10).Value = 2593.4443
$ws.Cells.Item(122, 11).Value = 7187.3079
$ws.Cells.Item(122, 12).Value = 7780.3329
$ws.Cells.Item(122, 13).Value = -4737.3079
$ws.Cells.Item(122, 14).Value = -12680.3329

$ws.Cells.Item(132, 8).Value = 3360.6
$ws.Cells.Item(132, 9).Value = 3414.2942
$ws.Cells.Item(132, 10).Value = 3165
$ws.Cells.Item(132, 11).Value = 10242.8826
$ws.Cells.Item(132, 12).Value = 9495
$ws.Cells.Item(132, 13).Value = -7712.882599999999
$ws.Cells.Item(132, 14).Value = -14555

$ws.Cells.Item(136, 8).Value = 1966.1212
$ws.Cells.Item(136, 9).Value = 1799.36
$ws.Cells.Item(136, 10).Value = 2487.25
$ws.Cells.Item(136, 11).Value = 5398.08
$ws.Cells.Item(136, 12).Value = 7461.75
$ws.Cells.Item(136, 13).Value = -2848.08
$ws.Cells.Item(136, 14).Value = -12561.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(50, 8).Value = 30780
$ws.Cells.Item(50, 10).Value = 30780
$ws.Cells.Item(50, 12).Value = 30780
$ws.Cells.Item(50, 14).Value = -31928

$ws.Cells.Item(107, 8).Value = 17052.441
$ws.Cells.Item(107, 9).Value = 19395.586
$ws.Cells.Item(107, 10).Value = 3462.2
$ws.Cells.Item(107, 11).Value = 19395.586
$ws.Cells.Item(107, 12).Value = 3462.2
$ws.Cells.Item(107, 13).Value = -17475.586
$ws.Cells.Item(107, 14).Value = -7302.2

$ws.Cells.Item(134, 8).Value = 2056.5356
$ws.Cells.Item(134, 9).Value = 1815.1818
$ws.Cells.Item(134, 10).Value = 2941.5
$ws.Cells.Item(134, 11).Value = 5445.5454
$ws.Cells.Item(134, 12).Value = 8824.5
$ws.Cells.Item(134, 13).Value = -2910.5454
$ws.Cells.Item(134, 14).Value = -13894.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2349.5151
$ws.Cells.Item(31, 9).Value = 1296.7
$ws.Cells.Item(31, 10).Value = 3969.2307
$ws.Cells.Item(31, 11).Value = 1296.7
$ws.Cells.Item(31, 12).Value = 3969.2307
$ws.Cells.Item(31, 13).Value = -1001.7
$ws.Cells.Item(31, 14).Value = -4559.2307

$ws.Cells.Item(34, 8).Value = 2349.5151
$ws.Cells.Item(34, 9).Value = 1296.7
$ws.Cells.Item(34, 10).Value = 3969.2307
$ws.Cells.Item(34, 11).Value = 1296.7
$ws.Cells.Item(34, 12).Value = 3969.2307
$ws.Cells.Item(34, 13).Value = -1094.7
$ws.Cells.Item(34, 14).Value = -4373.2307

$ws.Cells.Item(58, 8).Value = 700448.7
$ws.Cells.Item(58, 9).Value = 951118.4
$ws.Cells.Item(58, 10).Value = 2154.5715
$ws.Cells.Item(58, 11).Value = 951118.4
$ws.Cells.Item(58, 12).Value = 2154.5715
$ws.Cells.Item(58, 13).Value = -950915.4
$ws.Cells.Item(58, 14).Value = -2560.5715

$ws.Cells.Item(86, 8).Value = 3582.5454
$ws.Cells.Item(86, 9).Value = 3500
$ws.Cells.Item(86, 10).Value = 3600.889
$ws.Cells.Item(86, 11).Value = 3500
$ws.Cells.Item(86, 12).Value = 3600.889
$ws.Cells.Item(86, 13).Value = -2377
$ws.Cells.Item(86, 14).Value = -5846.889

$ws.Cells.Item(89, 8).Value = 3582.5454
$ws.Cells.Item(89, 9).Value = 3500
$ws.Cells.Item(89, 10).Value = 3600.889
$ws.Cells.Item(89, 11).Value = 17500
$ws.Cells.Item(89, 12).Value = 18004.445
$ws.Cells.Item(89, 13).Value = -11884
$ws.Cells.Item(89, 14).Value = -29236.445

$ws.Cells.Item(120, 8).Value = 29979
$ws.Cells.Item(120, 10).Value = 29979
$ws.Cells.Item(120, 12).Value = 29979
$ws.Cells.Item(120, 14).Value = -37237

$ws.Cells.Item(121, 8).Value = 29737.5
$ws.Cells.Item(121, 10).Value = 29737.5
$ws.Cells.Item(121, 12).Value = 29737.5
$ws.Cells.Item(121, 14).Value = -32357.5

$ws.Cells.Item(132, 8).Value = 423905.53
$ws.Cells.Item(132, 9).Value = 467292.38
$ws.Cells.Item(132, 10).Value = 4499.3335
$ws.Cells.Item(132, 11).Value = 1401877.14
$ws.Cells.Item(132, 12).Value = 13498.0005
$ws.Cells.Item(132, 13).Value = -1399347.14
$ws.Cells.Item(132, 14).Value = -18558.0005

$ws.Cells.Item(134, 8).Value = 1718.4546
$ws.Cells.Item(134, 9).Value = 1155.1111
$ws.Cells.Item(134, 10).Value = 4253.5
$ws.Cells.Item(134, 11).Value = 3465.3333
$ws.Cells.Item(134, 12).Value = 12760.5
$ws.Cells.Item(134, 13).Value = -930.3333000000002
$ws.Cells.Item(134, 14).Value = -17830.5

$ws.Cells.Item(136, 8).Value = 700448.7
$ws.Cells.Item(136, 9).Value = 951118.4
$ws.Cells.Item(136, 10).Value = 2154.5715
$ws.Cells.Item(136, 11).Value = 2853355.2
$ws.Cells.Item(136, 12).Value = 6463.7145
$ws.Cells.Item(136, 13).Value = -2850805.2
$ws.Cells.Item(136, 14).Value = -11563.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 11365135
$ws.Cells.Item(131, 9).Value = 1797.1428
$ws.Cells.Item(131, 10).Value = 12347152
$ws.Cells.Item(131, 11).Value = 5391.428400000001
$ws.Cells.Item(131, 12).Value = 37041456
$ws.Cells.Item(131, 13).Value = -351.4284000000007
$ws.Cells.Item(131, 14).Value = -37051536

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(45, 8).Value = 13459.5
$ws.Cells.Item(45, 10).Value = 13459.5
$ws.Cells.Item(45, 12).Value = 13459.5
$ws.Cells.Item(45, 14).Value = -14577.5

$ws.Cells.Item(51, 8).Value = 35559.6
$ws.Cells.Item(51, 10).Value = 35559.6
$ws.Cells.Item(51, 12).Value = 35559.6
$ws.Cells.Item(51, 14).Value = -36577.6

$ws.Cells.Item(132, 8).Value = 2219.6099
$ws.Cells.Item(132, 9).Value = 1617.2667
$ws.Cells.Item(132, 10).Value = 3862.3635
$ws.Cells.Item(132, 11).Value = 4851.800099999999
$ws.Cells.Item(132, 12).Value = 11587.0905
$ws.Cells.Item(132, 13).Value = -2321.800099999999
$ws.Cells.Item(132, 14).Value = -16647.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3663.6428
$ws.Cells.Item(132, 9).Value = 3300.2104
$ws.Cells.Item(132, 10).Value = 4430.8887
$ws.Cells.Item(132, 11).Value = 9900.6312
$ws.Cells.Item(132, 12).Value = 13292.6661
$ws.Cells.Item(132, 13).Value = -7370.6312
$ws.Cells.Item(132, 14).Value = -18352.6661

$ws.Cells.Item(136, 8).Value = 25252272
$ws.Cells.Item(136, 9).Value = 41668144
$ws.Cells.Item(136, 11).Value = 125004432
$ws.Cells.Item(136, 13).Value = -125001882

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 217549.83
$ws.Cells.Item(81, 9).Value = 250874.75
$ws.Cells.Item(81, 10).Value = 150900
$ws.Cells.Item(81, 11).Value = 501749.5
$ws.Cells.Item(81, 12).Value = 301800
$ws.Cells.Item(81, 13).Value = -500688.5
$ws.Cells.Item(81, 14).Value = -303922

$ws.Cells.Item(84, 8).Value = 217549.83
$ws.Cells.Item(84, 9).Value = 250874.75
$ws.Cells.Item(84, 10).Value = 150900
$ws.Cells.Item(84, 11).Value = 2508747.5
$ws.Cells.Item(84, 12).Value = 1509000
$ws.Cells.Item(84, 13).Value = -2503443.5
$ws.Cells.Item(84, 14).Value = -1519608

$ws.Cells.Item(132, 8).Value = 1576.2712
$ws.Cells.Item(132, 9).Value = 984.6667
$ws.Cells.Item(132, 10).Value = 2729.9
$ws.Cells.Item(132, 11).Value = 2954.0001
$ws.Cells.Item(132, 12).Value = 8189.700000000001
$ws.Cells.Item(132, 13).Value = -424.0001000000002
$ws.Cells.Item(132, 14).Value = -13249.7

$ws.Cells.Item(136, 8).Value = 1719.4062
$ws.Cells.Item(136, 9).Value = 1634
$ws.Cells.Item(136, 10).Value = 1975.625
$ws.Cells.Item(136, 11).Value = 4902
$ws.Cells.Item(136, 12).Value = 5926.875
$ws.Cells.Item(136, 13).Value = -2352
$ws.Cells.Item(136, 14).Value = -11026.875
